$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 11:17:22"
$wsZhCn.Range("H2").Value = "2016-03-23 11:17:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 11:17:26"
$wsDeDe.Range("H2").Value = "2016-03-23 11:17:55"
